$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @('D2', '29.305.09'),
    @('E2', '  +0.09%  '),
    @('D3', '1.841.60'),
    @('E3', '  -0.03%  '),
    @('D4', '0.9982'),
    @('E4', '  -0.18%  '),
    @('D5', '240.80'),
    @('E5', '  -1.38%  '),
    @('D6', '0.6692'),
    @('E6', '  -2.53%  '),
    @('D7', '0.9993'),
    @('E7', '  -0.12%  '),
    @('D8', '0.07422'),
    @('E8', '  -1.46%  '),
    @('D9', '0.2960'),
    @('E9', '  -2.06%  '),
    @('D10', '22.90'),
    @('E10', '  -1.46%  '),
    @('E11', '  +0.76%  '),
    @('D12', '5.036'),
    @('E12', '  -0.87%  '),
    @('D13', '1.790.38'),
    @('E13', '  -3.67%  '),
    @('D14', '0.6800'),
    @('E14', '  -0.71%  '),
    @('D15', '86.55'),
    @('E15', '  -3.02%  '),
    @('D16', '6.210'),
    @('E16', '  -1.50%  '),
    @('D17', '29.360.13'),
    @('E17', '  +0.25%  '),
    @('D18', '0.000008250'),
    @('E18', '  +0.41%  '),
    @('D19', '229.71'),
    @('E19', '  -2.03%  '),
    @('D20', '12.55'),
    @('E20', '  -0.18%  '),
    @('D21', '0.9986'),
    @('E21', '  -0.17%  '),
    @('D22', '7.296'),
    @('E22', '  -3.13%  '),
    @('D23', '0.9994'),
    @('E23', '  -0.11%  '),
    @('D24', '160.22'),
    @('E24', '  +0.27%  '),
    @('D25', '8.730'),
    @('E25', '  -1.07%  '),
    @('D26', '0.1416'),
    @('E26', '  -2.75%  '),
    @('D27', '18.05'),
    @('E27', '  -0.04%  '),
    @('D28', '1.509'),
    @('E28', '  -0.99%  '),
    @('D29', '4.213'),
    @('E29', '  -0.05%  '),
    @('D30', '4.090'),
    @('E30', '  -0.89%  '),
    @('D31', '1.200'),
    @('E31', '  +0.01%  '),
    @('D32', '0.05343'),
    @('E32', '  +3.77%  '),
    @('B33', 'ImmutableX'),
    @('C33', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D33', '0.7584'),
    @('E33', '  -1.62%  '),
    @('B34', 'LidoDAOToken'),
    @('C34', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'),
    @('D34', '1.867'),
    @('E34', '  +0.68%  '),
    @('D35', '1.137'),
    @('E35', '  -0.05%  '),
    @('D36', '2.680'),
    @('E36', '  +0.18%  '),
    @('D37', '1.338.10'),
    @('E37', '  +3.39%  '),
    @('D38', '0.01802'),
    @('E38', '  -2.31%  '),
    @('D39', '2.734'),
    @('E39', '  +1.11%  '),
    @('D40', '0.9213'),
    @('E40', '  -2.22%  '),
    @('D41', '6.004'),
    @('E41', '  +5.81%  '),
    @('E42', '  +0.08%  '),
    @('B43', 'Quant'),
    @('C43', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'),
    @('D43', '103.49'),
    @('E43', '  -1.80%  '),
    @('B44', 'XinFinNetwork'),
    @('C44', 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'),
    @('D44', '0.08045'),
    @('E44', '  +17.51%  '),
    @('E45', '  -1.31%  '),
    @('B46', 'RocketPoolETH'),
    @('C46', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'),
    @('D46', '1.962.91'),
    @('E46', '  -1.52%  '),
    @('B47', 'Mantle'),
    @('C47', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D47', '0.5157'),
    @('E47', '  -0.83%  '),
    @('B48', 'Aave'),
    @('C48', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D48', '63.93'),
    @('E48', '  +1.48%  '),
    @('B49', 'RenderToken'),
    @('C49', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D49', '1.766'),
    @('E49', '  -0.07%  '),
    @('D50', '9.293'),
    @('E50', '  -3.72%  '),
    @('D51', '0.05960'),
    @('E51', '  +0.49%  ')
)

foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $text = $edit[1]
    $escaped = $text.Replace('"', '""')
    $r = $ws.Range($cellRef)
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
